# Update res_bus vm_pu values for the 380 kV case (rows 2-25, columns B-F and I-N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.027102853657206
$ws.Cells.Item(2, 4).Value = 1.030826718843157
$ws.Cells.Item(2, 5).Value = 1.035877072077094
$ws.Cells.Item(2, 6).Value = 1.043634919974799
$ws.Cells.Item(2, 9).Value = 1.032995426893294
$ws.Cells.Item(2, 10).Value = 1.032262592691959
$ws.Cells.Item(2, 11).Value = 1.033636605255288
$ws.Cells.Item(2, 12).Value = 1.038672419361101
$ws.Cells.Item(2, 13).Value = 1.046408231202599
$ws.Cells.Item(2, 14).Value = 1.014769821801377

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027933215800363
$ws.Cells.Item(3, 4).Value = 1.031254391436562
$ws.Cells.Item(3, 5).Value = 1.036649267784965
$ws.Cells.Item(3, 6).Value = 1.044568126639018
$ws.Cells.Item(3, 9).Value = 1.033106422492641
$ws.Cells.Item(3, 10).Value = 1.032733609880367
$ws.Cells.Item(3, 11).Value = 1.033873736513128
$ws.Cells.Item(3, 12).Value = 1.039254199749573
$ws.Cells.Item(3, 13).Value = 1.047152193618798
$ws.Cells.Item(3, 14).Value = 1.014927455671785

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.02847115043483
$ws.Cells.Item(4, 4).Value = 1.031531439199948
$ws.Cells.Item(4, 5).Value = 1.03714994214315
$ws.Cells.Item(4, 6).Value = 1.045173344094577
$ws.Cells.Item(4, 9).Value = 1.033177245261059
$ws.Cells.Item(4, 10).Value = 1.033038371902111
$ws.Cells.Item(4, 11).Value = 1.034026757204192
$ws.Cells.Item(4, 12).Value = 1.03963100222943
$ws.Cells.Item(4, 13).Value = 1.04763432957096
$ws.Cells.Item(4, 14).Value = 1.015029404641446

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028697448542754
$ws.Cells.Item(5, 4).Value = 1.031647983615646
$ws.Cells.Item(5, 5).Value = 1.037360666174523
$ws.Cells.Item(5, 6).Value = 1.045428103748726
$ws.Cells.Item(5, 9).Value = 1.033206779346
$ws.Cells.Item(5, 10).Value = 1.033166488432122
$ws.Cells.Item(5, 11).Value = 1.034090985593166
$ws.Cells.Item(5, 12).Value = 1.039789492790809
$ws.Cells.Item(5, 13).Value = 1.047837195689898
$ws.Cells.Item(5, 14).Value = 1.015072251400398

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028735453763617
$ws.Cells.Item(6, 4).Value = 1.03166755619338
$ws.Cells.Item(6, 5).Value = 1.037396061720103
$ws.Cells.Item(6, 6).Value = 1.045470898093225
$ws.Cells.Item(6, 9).Value = 1.033211724166361
$ws.Cells.Item(6, 10).Value = 1.03318799940798
$ws.Cells.Item(6, 11).Value = 1.034101763835393
$ws.Cells.Item(6, 12).Value = 1.039816108868408
$ws.Cells.Item(6, 13).Value = 1.047871268087929
$ws.Cells.Item(6, 14).Value = 1.015079444808406

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.028474173652283
$ws.Cells.Item(7, 4).Value = 1.031532996185561
$ws.Cells.Item(7, 5).Value = 1.037152756905065
$ws.Cells.Item(7, 6).Value = 1.045176746925949
$ws.Cells.Item(7, 9).Value = 1.033177640840087
$ws.Cells.Item(7, 10).Value = 1.033040083824287
$ws.Cells.Item(7, 11).Value = 1.034027615827418
$ws.Cells.Item(7, 12).Value = 1.039633119664468
$ws.Cells.Item(7, 13).Value = 1.047637039587763
$ws.Cells.Item(7, 14).Value = 1.015029977211968

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.027383346036489
$ws.Cells.Item(8, 4).Value = 1.030971186450627
$ws.Cells.Item(8, 5).Value = 1.03613782860686
$ws.Cells.Item(8, 6).Value = 1.043950016698229
$ws.Cells.Item(8, 9).Value = 1.033033144769739
$ws.Cells.Item(8, 10).Value = 1.032421778070389
$ws.Cells.Item(8, 11).Value = 1.033716830941399
$ws.Cells.Item(8, 12).Value = 1.038868961115796
$ws.Cells.Item(8, 13).Value = 1.046659502079425
$ws.Cells.Item(8, 14).Value = 1.014823105127783

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.025466100460213
$ws.Cells.Item(9, 4).Value = 1.029983701147956
$ws.Cells.Item(9, 5).Value = 1.03435722493918
$ws.Cells.Item(9, 6).Value = 1.041798937367609
$ws.Cells.Item(9, 9).Value = 1.032770905137321
$ws.Cells.Item(9, 10).Value = 1.031332163175015
$ws.Cells.Item(9, 11).Value = 1.033166030584522
$ws.Cells.Item(9, 12).Value = 1.037525170989457
$ws.Cells.Item(9, 13).Value = 1.044942713587322
$ws.Cells.Item(9, 14).Value = 1.014458201953948

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.024191347268759
$ws.Cells.Item(10, 4).Value = 1.029327167238278
$ws.Cells.Item(10, 5).Value = 1.033175525578265
$ws.Cells.Item(10, 6).Value = 1.040372104281958
$ws.Cells.Item(10, 9).Value = 1.032590997386193
$ws.Cells.Item(10, 10).Value = 1.030605770949654
$ws.Cells.Item(10, 11).Value = 1.032796777706697
$ws.Cells.Item(10, 12).Value = 1.03663124611321
$ws.Cells.Item(10, 13).Value = 1.043802154229004
$ws.Cells.Item(10, 14).Value = 1.014214712711912

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.023640193858967
$ws.Cells.Item(11, 4).Value = 1.02904332905881
$ws.Cells.Item(11, 5).Value = 1.032665132638156
$ws.Cells.Item(11, 6).Value = 1.039756005248364
$ws.Cells.Item(11, 9).Value = 1.032511898853277
$ws.Cells.Item(11, 10).Value = 1.030291254535071
$ws.Cells.Item(11, 11).Value = 1.032636415201276
$ws.Cells.Item(11, 12).Value = 1.036244643289697
$ws.Cells.Item(11, 13).Value = 1.043309240289196
$ws.Cells.Item(11, 14).Value = 1.014109232598262

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.023435596200956
$ws.Cells.Item(12, 4).Value = 1.028937967770634
$ws.Cells.Item(12, 5).Value = 1.032475745693264
$ws.Cells.Item(12, 6).Value = 1.039527420276158
$ws.Cells.Item(12, 9).Value = 1.032482338950316
$ws.Cells.Item(12, 10).Value = 1.030174432794509
$ws.Cells.Item(12, 11).Value = 1.032576779486308
$ws.Cells.Item(12, 12).Value = 1.036101114194213
$ws.Cells.Item(12, 13).Value = 1.043126295323619
$ws.Cells.Item(12, 14).Value = 1.014070045885065

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.023479477399314
$ws.Cells.Item(13, 4).Value = 1.028960564977908
$ws.Cells.Item(13, 5).Value = 1.032516360950138
$ws.Cells.Item(13, 6).Value = 1.039576440658937
$ws.Cells.Item(13, 9).Value = 1.032488687747252
$ws.Cells.Item(13, 10).Value = 1.030199491269155
$ws.Cells.Item(13, 11).Value = 1.032589574695185
$ws.Cells.Item(13, 12).Value = 1.036131898373186
$ws.Cells.Item(13, 13).Value = 1.043165531042243
$ws.Cells.Item(13, 14).Value = 1.014078451864009

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.023623279188276
$ws.Cells.Item(14, 4).Value = 1.029034618446145
$ws.Cells.Item(14, 5).Value = 1.032649473851161
$ws.Cells.Item(14, 6).Value = 1.039737104995087
$ws.Cells.Item(14, 9).Value = 1.032509459073805
$ws.Cells.Item(14, 10).Value = 1.030281597931816
$ws.Cells.Item(14, 11).Value = 1.032631487114257
$ws.Cells.Item(14, 12).Value = 1.036232777644648
$ws.Cells.Item(14, 13).Value = 1.043294115029467
$ws.Cells.Item(14, 14).Value = 1.014105993546402

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023711896862164
$ws.Cells.Item(15, 4).Value = 1.029080254413204
$ws.Cells.Item(15, 5).Value = 1.032731515091418
$ws.Cells.Item(15, 6).Value = 1.039836130320841
$ws.Cells.Item(15, 9).Value = 1.0325222332505
$ws.Cells.Item(15, 10).Value = 1.030332187078106
$ws.Cells.Item(15, 11).Value = 1.032657301507251
$ws.Cells.Item(15, 12).Value = 1.036294942330061
$ws.Cells.Item(15, 13).Value = 1.04337335915924
$ws.Cells.Item(15, 14).Value = 1.014122962011486

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.024227942983669
$ws.Cells.Item(16, 4).Value = 1.029346014192201
$ws.Cells.Item(16, 5).Value = 1.033209426024266
$ws.Cells.Item(16, 6).Value = 1.040413029334923
$ws.Cells.Item(16, 9).Value = 1.032596221728965
$ws.Cells.Item(16, 10).Value = 1.030626644814289
$ws.Cells.Item(16, 11).Value = 1.032807410563658
$ws.Cells.Item(16, 12).Value = 1.036656913745875
$ws.Cells.Item(16, 13).Value = 1.043834887579614
$ws.Cells.Item(16, 14).Value = 1.014221712114572

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024551866580868
$ws.Cells.Item(17, 4).Value = 1.029512838883501
$ws.Cells.Item(17, 5).Value = 1.033509553754586
$ws.Cells.Item(17, 6).Value = 1.040775367089815
$ws.Cells.Item(17, 9).Value = 1.032642312687555
$ws.Cells.Item(17, 10).Value = 1.030811355577047
$ws.Cells.Item(17, 11).Value = 1.032901444074389
$ws.Cells.Item(17, 12).Value = 1.036884096362499
$ws.Cells.Item(17, 13).Value = 1.044124649200691
$ws.Cells.Item(17, 14).Value = 1.014283643007852

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024740885048307
$ws.Cells.Item(18, 4).Value = 1.029610187742089
$ws.Cells.Item(18, 5).Value = 1.033684737592688
$ws.Cells.Item(18, 6).Value = 1.040986879280845
$ws.Cells.Item(18, 9).Value = 1.032669081196308
$ws.Cells.Item(18, 10).Value = 1.030919095722051
$ws.Cells.Item(18, 11).Value = 1.032956246456522
$ws.Cells.Item(18, 12).Value = 1.037016653638348
$ws.Cells.Item(18, 13).Value = 1.04429375436276
$ws.Cells.Item(18, 14).Value = 1.01431976162038

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.02480534886561
$ws.Cells.Item(19, 4).Value = 1.029643388394089
$ws.Cells.Item(19, 5).Value = 1.033744491801543
$ws.Cells.Item(19, 6).Value = 1.041059027692137
$ws.Cells.Item(19, 9).Value = 1.032678188942682
$ws.Cells.Item(19, 10).Value = 1.030955832541416
$ws.Cells.Item(19, 11).Value = 1.032974924851764
$ws.Cells.Item(19, 12).Value = 1.037061859918339
$ws.Cells.Item(19, 13).Value = 1.044351430427236
$ws.Cells.Item(19, 14).Value = 1.014332076334624

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.024517104429211
$ws.Cells.Item(20, 4).Value = 1.029494935727921
$ws.Cells.Item(20, 5).Value = 1.03347734000807
$ws.Cells.Item(20, 6).Value = 1.040736474420933
$ws.Cells.Item(20, 9).Value = 1.0326373795126
$ws.Cells.Item(20, 10).Value = 1.030791537698151
$ws.Cells.Item(20, 11).Value = 1.032891359894487
$ws.Cells.Item(20, 12).Value = 1.036859717101388
$ws.Cells.Item(20, 13).Value = 1.044093550981862
$ws.Cells.Item(20, 14).Value = 1.014276998884676

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.023580929675042
$ws.Cells.Item(21, 4).Value = 1.029012809629337
$ws.Cells.Item(21, 5).Value = 1.032610270002259
$ws.Cells.Item(21, 6).Value = 1.039689786115018
$ws.Cells.Item(21, 9).Value = 1.032503347376687
$ws.Cells.Item(21, 10).Value = 1.030257419455261
$ws.Cells.Item(21, 11).Value = 1.032619146872487
$ws.Cells.Item(21, 12).Value = 1.036203069203433
$ws.Cells.Item(21, 13).Value = 1.043256246228319
$ws.Cells.Item(21, 14).Value = 1.014097883388735

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022993044948302
$ws.Cells.Item(22, 4).Value = 1.028710077113642
$ws.Cells.Item(22, 5).Value = 1.032066241936422
$ws.Cells.Item(22, 6).Value = 1.039033206329809
$ws.Cells.Item(22, 9).Value = 1.032418039681894
$ws.Cells.Item(22, 10).Value = 1.029921619896871
$ws.Cells.Item(22, 11).Value = 1.032447591684008
$ws.Cells.Item(22, 12).Value = 1.035790628229723
$ws.Cells.Item(22, 13).Value = 1.042730640491966
$ws.Cells.Item(22, 14).Value = 1.01398522781572

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.023304624463684
$ws.Cells.Item(23, 4).Value = 1.028870522852712
$ws.Cells.Item(23, 5).Value = 1.032354533483733
$ws.Cells.Item(23, 6).Value = 1.03938112749889
$ws.Cells.Item(23, 9).Value = 1.032463360909366
$ws.Cells.Item(23, 10).Value = 1.03009963112295
$ws.Cells.Item(23, 11).Value = 1.032538574232701
$ws.Cells.Item(23, 12).Value = 1.03600923074658
$ws.Cells.Item(23, 13).Value = 1.043009193778853
$ws.Cells.Item(23, 14).Value = 1.014044952167649

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024532811704348
$ws.Cells.Item(24, 4).Value = 1.029503025261196
$ws.Cells.Item(24, 5).Value = 1.033491895629317
$ws.Cells.Item(24, 6).Value = 1.040754047829013
$ws.Cells.Item(24, 9).Value = 1.03263960895919
$ws.Cells.Item(24, 10).Value = 1.030800492540144
$ws.Cells.Item(24, 11).Value = 1.032895916642911
$ws.Cells.Item(24, 12).Value = 1.036870732899681
$ws.Cells.Item(24, 13).Value = 1.044107602643936
$ws.Cells.Item(24, 14).Value = 1.014280001092309

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025961160201758
$ws.Cells.Item(25, 4).Value = 1.030238682504893
$ws.Cells.Item(25, 5).Value = 1.03481661493955
$ws.Cells.Item(25, 6).Value = 1.042353778390405
$ws.Cells.Item(25, 9).Value = 1.032839599015142
$ws.Cells.Item(25, 10).Value = 1.031613856917101
$ws.Cells.Item(25, 11).Value = 1.033308792614253
$ws.Cells.Item(25, 12).Value = 1.03787223793333
$ws.Cells.Item(25, 13).Value = 1.045385852459704
$ws.Cells.Item(25, 14).Value = 1.014552579008443

Write-Host "Updated vm_pu values for rows 2-25 (380 kV case)"
